# Updated symbol list on Mon Jan  2 17:32:36 UTC 2023 with GitHub Actions
# This script updates Price (column D) and Volume(1h) (column E) cells
# in the active worksheet, preserving their original text/string cell type
# (values in this sheet are stored as text, e.g. "246.53" or "0.63%", not
# numbers, so we force a Text number format while writing the value and
# then restore the cell's original style to avoid changing formatting).

function Set-CellText($ws, $ref, $val) {
    $cell = $ws.Range($ref)
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = $origStyle
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

Set-CellText $ws "D2" "246.53"
Set-CellText $ws "E2" "0.63%"
Set-CellText $ws "D3" "29.56"
Set-CellText $ws "E3" "9.73%"
Set-CellText $ws "D4" "5.160"
Set-CellText $ws "E4" "1.65%"
Set-CellText $ws "D5" "0.05707"
Set-CellText $ws "D6" "6.600"
Set-CellText $ws "E6" "1.86%"
Set-CellText $ws "D7" "3.073"
Set-CellText $ws "E7" "2.32%"
Set-CellText $ws "D8" "0.8578"
Set-CellText $ws "E8" "4.59%"
Set-CellText $ws "D9" "0.8690"
Set-CellText $ws "E9" "3.64%"
Set-CellText $ws "D10" "0.1365"
Set-CellText $ws "E10" "2.79%"
Set-CellText $ws "D11" "0.07085"
Set-CellText $ws "E11" "2.62%"
Set-CellText $ws "D12" "0.02924"
Set-CellText $ws "E12" "2.30%"
Set-CellText $ws "D13" "0.09381"
Set-CellText $ws "E13" "-0.22%"
Set-CellText $ws "D14" "0.001522"
Set-CellText $ws "E14" "0.10%"
Set-CellText $ws "D15" "0.04182"
Set-CellText $ws "D16" "0.0006010"
Set-CellText $ws "E16" "0.30%"
Set-CellText $ws "D17" "0.006158"
Set-CellText $ws "E17" "1.05%"
Set-CellText $ws "D19" "3.485"
Set-CellText $ws "E19" "-0.72%"
Set-CellText $ws "D20" "2.181"
Set-CellText $ws "E20" "-5.80%"
Set-CellText $ws "D21" "0.3174"
Set-CellText $ws "E21" "-0.10%"
Set-CellText $ws "D22" "0.03309"
Set-CellText $ws "E22" "4.21%"
Set-CellText $ws "D23" "0.1326"
Set-CellText $ws "E23" "2.18%"
Set-CellText $ws "D24" "3.481"
Set-CellText $ws "E24" "-1.91%"
Set-CellText $ws "E25" "0.55%"
Set-CellText $ws "D26" "0.005032"
Set-CellText $ws "E26" "26.85%"
Set-CellText $ws "E27" "0.12%"
Set-CellText $ws "E28" "23.55%"
Set-CellText $ws "D40" "0.03744"
Set-CellText $ws "E40" "1.33%"
Set-CellText $ws "D41" "0.005769"
Set-CellText $ws "E41" "-3.13%"
Set-CellText $ws "D42" "0.1072"
Set-CellText $ws "E42" "1.48%"
Set-CellText $ws "D43" "0.002540"
Set-CellText $ws "E43" "10.50%"
Set-CellText $ws "D44" "0.009964"
Set-CellText $ws "E44" "6.08%"
Set-CellText $ws "E45" "0.25%"
Set-CellText $ws "E46" "0.08%"
Set-CellText $ws "D47" "0.06000"
Set-CellText $ws "D48" "0.002563"
Set-CellText $ws "E48" "-1.18%"
Set-CellText $ws "E49" "0.08%"
Set-CellText $ws "E50" "0.08%"
